$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "38.031.09"
$ws.Range("E2").Value = "  +2.59%  "

# Row 3
$ws.Range("D3").Value = "2.048.65"
$ws.Range("E3").Value = "  +1.53%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.43%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.615"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.87%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.49"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.40%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.384"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.02%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0801"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.86%  "

# Row 11
$ws.Range("E11").Value = "  +1.76%  "

# Row 12
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "2.356.65"
$ws.Range("E12").Value = "  +1.73%  "

# Row 13
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.54%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.15%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.96%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.755"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.49%  "

# Row 17
$ws.Range("D17").Value = "2.049.39"
$ws.Range("E17").Value = "  +1.10%  "

# Row 18
$ws.Range("D18").Value = "37.993.77"
$ws.Range("E18").Value = "  +2.61%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.88%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.16%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0827"
$ws.Range("E21").Value = "  +1.55%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "225.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.03%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.04%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.44%  "

# Row 25
$ws.Range("E25").Value = "  +1.54%  "

# Row 26
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.05%  "

# Row 27
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.11%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.132"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.85%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.97"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.59%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.44%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.119"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.94%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.50"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.66%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.55"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.17%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.04"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.88%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0602"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.06%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.24"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +14.34%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.36%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.26"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.25%  "

# Row 39
$ws.Range("E39").Value = "  +0.01%  "

# Row 40
$ws.Range("D40").Value = "1.518.23"
$ws.Range("E40").Value = "  +3.32%  "

# Row 41
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "96.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.53%  "

# Row 42
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0215"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.37%  "

# Row 43
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.60%  "

# Row 44
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.59%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0923"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.48%  "

# Row 46
$ws.Range("E46").Value = "  +1.77%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.65%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.59%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.41%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.03"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.19%  "

# Row 51
$ws.Range("D51").Value = "2.245.10"
$ws.Range("E51").Value = "  +1.80%  "
